$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1Precios base SQL")

# --- Update rows 63-73: unify unit to "varilla" and normalize Fierro names to uppercase city suffix ---
$ws.Range("B63").Value = "varilla"
$ws.Range("C63").Value = 'Fierro 1/2" - AREQUIPA'

$ws.Range("B64").Value = "varilla"
$ws.Range("C64").Value = 'Fierro 1/2" - SIDERPERU'

$ws.Range("B65").Value = "varilla"
$ws.Range("C65").Value = 'Fierro 1/4" 6mm - AREQUIPA'

$ws.Range("B66").Value = "varilla"
$ws.Range("C66").Value = 'Fierro 1/4" 6mm - SIDERPERU'

$ws.Range("B67").Value = "varilla"
$ws.Range("C67").Value = 'Fierro 3/4" - SIDERPERU'

$ws.Range("B68").Value = "varilla"
$ws.Range("C68").Value = 'Fierro 3/8" - AREQUIPA'

$ws.Range("B69").Value = "varilla"
$ws.Range("C69").Value = 'Fierro 3/8" - SIDERPERU'

$ws.Range("B70").Value = "varilla"
$ws.Range("C70").Value = 'Fierro 5/8" - AREQUIPA'

$ws.Range("B71").Value = "varilla"
$ws.Range("C71").Value = 'Fierro 5/8" - SIDERPERU'

$ws.Range("B72").Value = "varilla"
$ws.Range("C72").Value = 'Fierro 8mm - AREQUIPA'

$ws.Range("B73").Value = "varilla"
$ws.Range("C73").Value = 'Fierro 8mm - SIDERPERU'

# --- View state: scroll to show row 111, select N2:N156 ---
$ws.Activate()
$ws.Range("A111").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 111
$ws.Range("N2:N156").Select() | Out-Null

# --- Page setup: Letter size paper ---
$ws.PageSetup.PaperSize = 9

# --- "Maestra de unidades" sheet: move selection to B15 ---
$ws2 = $wb.Worksheets.Item("Maestra de unidades")
$ws2.Activate()
$ws2.Range("B15").Select() | Out-Null

# Reactivate the main sheet at the end (tabSelected)
$ws.Activate()
